$wb = $excel.ActiveWorkbook

# --- Sheet2 (IIQScenario1): insert a new row 4 with the PIN/instructions text ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

$ws2.Rows.Item(4).Insert()

$newText = "This questionnaire should take you approximately 30 minutes to complete. It does not need to be completed all at once.`nWe encourage you to take a break and return to complete the survey later, if needed.`nYou can save your answers by clicking the ""Save and Next"" button in the survey and closing your browser window.`nWhen you return to the survey, you will be asked for a PIN code, provided during the log in process.`nThis unique PIN code returns you to your previous spot in the questionnaire.`nThe study team requests that you complete the questionnaire within two months from the date you start the questionnaire."

$ws2.Range("A4").Value = $newText
$ws2.Range("B4").Value = $newText
$ws2.Rows.Item(4).RowHeight = 187

# --- Sheet1 (screenerScenario1): update view/selection ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A60").Select()
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B36").Select()

# --- Sheet3 (RASSurveyScenario1): update view/selection ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("A4").Select()

# --- Sheet2: finalize view/selection and make it the active tab ---
$ws2.Activate()
$ws2.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("A4").Select()

Write-Host "done"
